$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# F3 was a hard-coded input of 2200; set it to 0. Every other touched cell
# (F4:F66, M3:M66, N3:N66, O3:O66, and the two chart numCaches that read
# Sheet1!N and Sheet1!O) is a formula that depends on F3, so Excel's normal
# recalculation reproduces the rest of the diff automatically.
$ws.Range("F3").Value = 0

# Move the selection to match the author's final cursor position.
$ws.Range("F7").Select()
